$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 7
$ws.Range("E7").Value = 8

# Row 18
$ws.Range("E18").Value = 111

# Row 19
$ws.Range("E19").Value = 57

# Row 26
$ws.Range("E26").Value = 29
$ws.Range("F26").Value = 15
$ws.Range("H26").Value = 25

# Row 36
$ws.Range("E36").Value = 101

# Row 37
$ws.Range("E37").Value = 55
$ws.Range("F37").Value = 31
$ws.Range("H37").Value = 43

# Row 87
$ws.Range("F87").Value = 4
$ws.Range("H87").Value = 11

# Row 88
$ws.Range("E88").Value = 24

# Row 89
$ws.Range("E89").Value = 41
